$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "isEmpty" column (C) with header + boolean values
$ws.Range("C1").Value = "isEmpty"
$ws.Range("C2").Value = $true
$ws.Range("C3").Value = $false
$ws.Range("C4").Value = $true

# Match the formatting used by the adjacent column (header style / data style)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

$ws.Range("B2:B4").Copy()
$ws.Range("C2:C4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the active selection, matching the post-edit workbook state
$ws.Range("D5").Select()
